$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("21-09-2021", 114.55, 187.64, 101.05),
    @("22-09-2021", 114.24, 186.96, 100.8),
    @("23-09-2021", 114.38, 186.91, 100.96),
    @("24-09-2021", 114.74, 187.51, 101.28),
    @("27-09-2021", 114.81, 187.63, 101.31),
    @("28-09-2021", 115.68, 189.08, 102.06),
    @("29-09-2021", 115.85, 189.5, 102.15),
    @("30-09-2021", 116.28, 190.03, 102.47),
    @("01-10-2021", 117.6, 192.04, 103.66)
)

$startRow = 182
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $dateCell = $ws.Cells.Item($r, 1)
    # Enter the date string as a text formula first (so Excel does not
    # auto-convert the dd-mm-yyyy text into a date serial number), then
    # convert it to a plain static value via copy / paste-special so the
    # stored cell stays a normal shared-string text cell (no formula,
    # no special number format / style left behind).
    $dateCell.Formula = "=""" + $row[0] + """"
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163)

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$excel.CutCopyMode = 0
